$wb = $excel.ActiveWorkbook

# Rename "Sheet1" to "Research"
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Name = "Research"
